# Lab 5 submission: tweak the abstract / footnote paragraph styles.
#  1. Add a new "Abstract Title" paragraph style (based on Normal, followed
#     by the Abstract style).
#  2. Pull the "Abstract" style's space-before in from 15pt (300 twips) to
#     5pt (100 twips) now that AbstractTitle supplies the space above it.
#  3. Add a new "Footnote Block Text" paragraph style (based on Footnote
#     Text), mirroring the existing Block Text style but for footnotes.

$d = $word.ActiveDocument

# --- 1. New "Abstract Title" style --------------------------------------
$abstractTitle = $d.Styles.Add("AbstractTitle", 1)
$abstractTitle.NameLocal = "Abstract Title"
$abstractTitle.BaseStyle = $d.Styles("Normal")
$abstractTitle.NextParagraphStyle = $d.Styles("Abstract")
$abstractTitle.QuickStyle = $true

$abstractTitle.ParagraphFormat.KeepWithNext = $true
$abstractTitle.ParagraphFormat.KeepTogether = $true
$abstractTitle.ParagraphFormat.Alignment = 1
$abstractTitle.ParagraphFormat.SpaceBefore = 15
$abstractTitle.ParagraphFormat.SpaceAfter = 0

$abstractTitle.Font.Size = 10
$abstractTitle.Font.SizeBi = 10
$abstractTitle.Font.Bold = $true
$abstractTitle.Font.Color = 9067060

# --- 2. "Abstract" style: tighten the space above the paragraph ---------
$abstract = $d.Styles("Abstract")
$abstract.ParagraphFormat.SpaceBefore = 5

# --- 3. New "Footnote Block Text" style ----------------------------------
$footnoteBlockText = $d.Styles.Add("FootnoteBlockText", 1)
$footnoteBlockText.NameLocal = "Footnote Block Text"
$footnoteBlockText.BaseStyle = $d.Styles("Footnote Text")
$footnoteBlockText.NextParagraphStyle = $d.Styles("Footnote Text")
$footnoteBlockText.Priority = 9
$footnoteBlockText.UnhideWhenUsed = $true
$footnoteBlockText.QuickStyle = $true

$footnoteBlockText.ParagraphFormat.SpaceBefore = 5
$footnoteBlockText.ParagraphFormat.SpaceAfter = 5
$footnoteBlockText.ParagraphFormat.FirstLineIndent = 0
$footnoteBlockText.ParagraphFormat.LeftIndent = 24
$footnoteBlockText.ParagraphFormat.RightIndent = 24
